# Target change: the document's last two (empty) paragraphs become
#   - an empty paragraph
#   - a paragraph containing "Another edit." followed by the _GoBack
#     bookmark (moved here from the "But can I Push again..." paragraph,
#     where it used to sit).
$d = $word.ActiveDocument

# The very last paragraph in the document is currently empty; type the new
# sentence into it. A temporary "|" marker is appended right after the
# sentence so the inserted text is immediately followed by further run
# content -- that keeps the subsequent zero-length range (used to drop the
# bookmark right after "Another edit.") well anchored. The marker is
# deleted again once the bookmark is in place.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("Another edit.|")

# Locate the marker and collapse to a point immediately before it -- i.e.
# immediately after "Another edit." -- which is where _GoBack belongs now.
$markerRange = $d.Content
$markerRange.Find.Execute("|", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newBookmarkSpot = $markerRange.Duplicate
$newBookmarkSpot.Collapse(1)

# Relocate _GoBack: drop it from its old spot (end of "But can I Push
# again without a pull?") and recreate it at the new edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $newBookmarkSpot)

# Clean up the temporary marker character.
$markerRange.Delete()
